$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks first (targets will be re-added after data is in place)
$ws.Cells.Hyperlinks.Delete()

# Write data rows 2-13 (A:H)
$ws.Range("A2").Value2 = '2025-10-27 12:39:32'
$ws.Range("B2").Value2 = 'Amazon購入履歴の明細PDFを自動ダウンロード&自動リネームするシステム開発'
$ws.Range("C2").Value2 = 'システム開発'
$ws.Range("D2").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E2").Value2 = '期限情報なし'
$ws.Range("F2").Value2 = 'https://www.lancers.jp/work/detail/5421083'
$ws.Range("G2").Value2 = 113
$ws.Range("H2").Value2 = '◆開発,システム開発'

$ws.Range("A3").Value2 = '2025-10-27 12:39:32'
$ws.Range("B3").Value2 = '【RPA構築依頼】不動産問い合わせ対応自動化(アシロボ使用、報酬10万円)'
$ws.Range("C3").Value2 = 'システム開発'
$ws.Range("D3").Value2 = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value2 = '期限情報なし'
$ws.Range("F3").Value2 = 'https://www.lancers.jp/work/detail/5421443'
$ws.Range("G3").Value2 = 88
$ws.Range("H3").Value2 = '◆自動化'

$ws.Range("A4").Value2 = '2025-10-27 12:39:32'
$ws.Range("B4").Value2 = '初回 楽天RMSの配布型クーポンの自動登録システムの開発'
$ws.Range("C4").Value2 = 'システム開発'
$ws.Range("D4").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value2 = '期限情報なし'
$ws.Range("F4").Value2 = 'https://www.lancers.jp/work/detail/5421265'
$ws.Range("G4").Value2 = 78
$ws.Range("H4").Value2 = '◆開発'

$ws.Range("A5").Value2 = '2025-10-27 12:39:32'
$ws.Range("B5").Value2 = '【カフェ情報プラットフォーム開発】基本設計からリリースまで'
$ws.Range("C5").Value2 = 'システム開発'
$ws.Range("D5").Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value2 = '期限情報なし'
$ws.Range("F5").Value2 = 'https://www.lancers.jp/work/detail/5420868'
$ws.Range("G5").Value2 = 75
$ws.Range("H5").Value2 = '◆開発'

$ws.Range("A6").Value2 = '2025-10-27 12:39:32'
$ws.Range("B6").Value2 = '【簡単RPA構築】特定ツールからのデータ取得・Excel処理・スプレッドシートへの貼付け'
$ws.Range("C6").Value2 = 'システム開発'
$ws.Range("D6").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value2 = '期限情報なし'
$ws.Range("F6").Value2 = 'https://www.lancers.jp/work/detail/5421445'
$ws.Range("G6").Value2 = 68
$ws.Range("H6").Value2 = '◆ツール'

$ws.Range("A7").Value2 = '2025-10-27 12:39:32'
$ws.Range("B7").Value2 = '音声デシベル検知器の開発を手伝ってくれる方募集!'
$ws.Range("C7").Value2 = 'システム開発'
$ws.Range("D7").Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E7").Value2 = '期限情報なし'
$ws.Range("F7").Value2 = 'https://www.lancers.jp/work/detail/5421105'
$ws.Range("G7").Value2 = 68
$ws.Range("H7").Value2 = '◆開発'

$ws.Range("A8").Value2 = '2025-10-27 12:39:32'
$ws.Range("B8").Value2 = 'WordPressサイトのリニューアル作業依頼'
$ws.Range("C8").Value2 = 'システム開発'
$ws.Range("D8").Value2 = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E8").Value2 = '期限情報なし'
$ws.Range("F8").Value2 = 'https://www.lancers.jp/work/detail/5420971'
$ws.Range("G8").Value2 = 50
$ws.Range("H8").Value2 = '◇サイト ○WordPress'

$ws.Range("A9").Value2 = '2025-10-27 12:39:32'
$ws.Range("B9").Value2 = '【カンタン作業】サイト環境立ち上げ検証の作業!'
$ws.Range("C9").Value2 = 'システム開発'
$ws.Range("D9").Value2 = '~ 5,000 円 / 固定'
$ws.Range("E9").Value2 = '期限情報なし'
$ws.Range("F9").Value2 = 'https://www.lancers.jp/work/detail/5421230'
$ws.Range("G9").Value2 = 30
$ws.Range("H9").Value2 = '◇サイト'

$ws.Range("A10").Value2 = '2025-10-27 12:39:32'
$ws.Range("B10").Value2 = '〖リモート可〗Delphiエンジニア募集'
$ws.Range("C10").Value2 = 'システム開発'
$ws.Range("D10").Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value2 = '期限情報なし'
$ws.Range("F10").Value2 = 'https://www.lancers.jp/work/detail/5341051'
$ws.Range("G10").Value2 = 25
$ws.Range("H10").ClearContents()

$ws.Range("A11").Value2 = '2025-10-27 12:39:32'
$ws.Range("B11").Value2 = '【急募】Shopifyでのフォーム一体型LPコード作成依頼'
$ws.Range("C11").Value2 = 'システム開発'
$ws.Range("D11").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E11").Value2 = '期限情報なし'
$ws.Range("F11").Value2 = 'https://www.lancers.jp/work/detail/5421564'
$ws.Range("G11").Value2 = 13
$ws.Range("H11").ClearContents()

$ws.Range("A12").Value2 = '2025-10-27 12:39:32'
$ws.Range("B12").Value2 = '【急募】年末調整業務のマクロ作成依頼'
$ws.Range("C12").Value2 = 'システム開発'
$ws.Range("D12").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E12").Value2 = '期限情報なし'
$ws.Range("F12").Value2 = 'https://www.lancers.jp/work/detail/5421418'
$ws.Range("G12").Value2 = 13
$ws.Range("H12").ClearContents()

$ws.Range("A13").Value2 = '2025-10-27 12:39:32'
$ws.Range("B13").Value2 = '【10,000円1万枚】指定したURL先のHPのキャプチャー画像を作成お願い致します。'
$ws.Range("C13").Value2 = 'システム開発'
$ws.Range("D13").Value2 = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E13").Value2 = '期限情報なし'
$ws.Range("F13").Value2 = 'https://www.lancers.jp/work/detail/5421177'
$ws.Range("G13").Value2 = 10
$ws.Range("H13").ClearContents()

# Re-add hyperlinks for F2:F13 pointing at the URL already in each cell
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5421083') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5421443') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5421265') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5420868') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5421445') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5421105') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5420971') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5421230') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5341051') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5421564') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5421418') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5421177') | Out-Null

# Column B width: target raw OOXML width=46 -> ColumnWidth = 46 - 5/6
$ws.Columns.Item(2).ColumnWidth = 45.166666666666664

"done"